$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column = 想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 32
$ws1.Range("F5").Value  = 91
$ws1.Range("F6").Value  = 1393
$ws1.Range("F7").Value  = 148
$ws1.Range("F9").Value  = 32
$ws1.Range("F10").Value = 9448
$ws1.Range("F11").Value = 155
$ws1.Range("F12").Value = 101
$ws1.Range("F13").Value = 225
$ws1.Range("F14").Value = 180
$ws1.Range("F16").Value = 6442
$ws1.Range("F17").Value = 1068
$ws1.Range("F18").Value = 102
$ws1.Range("F19").Value = 46
$ws1.Range("F20").Value = 145

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 39

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 32
$ws4.Range("F5").Value  = 91
$ws4.Range("F6").Value  = 1393
$ws4.Range("F7").Value  = 148
$ws4.Range("F9").Value  = 32
$ws4.Range("F10").Value = 39
$ws4.Range("F12").Value = 9448
$ws4.Range("F13").Value = 155
$ws4.Range("F14").Value = 101
$ws4.Range("F15").Value = 225
$ws4.Range("F16").Value = 180
$ws4.Range("F18").Value = 6442
$ws4.Range("F19").Value = 1068
$ws4.Range("F20").Value = 102
$ws4.Range("F21").Value = 46
$ws4.Range("F22").Value = 145
